{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies a set of exact text replacements inside existing paragraph runs,\n// matching the translated-string update described in the diff.\n// Each entry's `old` text is searched for verbatim (case-sensitive, exact\n// match) in the document body and replaced in-place with `new`, preserving\n// the paragraph/run formatting already present around the text.\nconst replacements = [\n  {\n    \"old\": \"Amava akho nenkqubo ye-ParentText abalulekile kuphononongo lwethu. Singathanda ukuva malunga namava akho ngayo, kwaye yintoni mhlawumbi engakhange ikusebenzele wena. This interview is part of a study carried out by researchers from the Universities of Cape Town in South Africa and the University of Oxford in the United Kingdom. \",\n    \"new\": \"Amava akho nenkqubo ye-ParentText abalulekile kuphononongo lwethu. Singathanda ukuva malunga namava akho ngayo, kwaye yintoni mhlawumbi engakhange ikusebenzele wena. Oludliwano-ndlebe luyinxalenye yophononongo olwenziwa ngabaphandi abaphuma kwiDyunivesithi yaseKapa eMzantsi Afrika kunye neDyunivesithi yaseOxford eUnited Kingdom. \"\n  },\n  {\n    \"old\": \"Before you decide if you\u2019d like to be interviewed, it\u2019s important for you to know why we\u2019re doing this research and what participating in it would involve. All the information you might need is explained below, but if you have any questions about your participation or our study, please email the study team at swift@globalparenting.org or message us on WhatsApp at +27 XX XXX XXXX. Silapha ukuzokunceda wena!\",\n    \"new\": \"Ngaphambi kokuba ugqibe ekubeni ungathanda na ukuba nodliwano-ndlebe, kubalulekile ukuba wazi kutheni sisenza oluphando nje kwaye ukuthatha inxalenye kungaquka ntoni. Lonke ulwazi ozakuludinga luchaziwe ngezantsi, kodwa ukuba unayo nayiphi na imibuzo malunga nokuthatha inxaxheba okanye ngophando lwethu, ndicela u-imeyilele iqela lophononongo ku swift@globalparenting.org okanye uthumele umyalezo kuthi ku WhatsApp ku +27 XX XXX XXXX. Silapha ukuzokunceda wena!\"\n  },\n  {\n    \"old\": \"We would like to have a telephonic conversation with you which will last a maximum of 15 minutes. Omnye wabaphandi bethu uzakutsalela umnxeba athethe nawe ngexesha elikulungeleyo wena. Akukho zimpendulu zilungileyo okanye ezingalunganga, sifuna nje ukuva amava kunye nemibono yakho nge chatbot. Please make sure that when we call, that you only let the interview start when you are in a private space where you feel comfortable to talk without being overheard or interrupted. Ukuba ngelixa wenziwa udliwano-ndlebe, uye waphazamiseka, ndicela ucele umphandi ukuba ame ude uzive ukhuselekile ukuqhubeka nokuthetha.\",\n    \"new\": \"Singathanda ukuba nencoko ngomnxeba nawe ozakuthatha imizuzu eyi-15 ubude. Omnye wabaphandi bethu uzakutsalela umnxeba athethe nawe ngexesha elikulungeleyo wena. Akukho zimpendulu zilungileyo okanye ezingalunganga, sifuna nje ukuva amava kunye nemibono yakho nge chatbot. Nceda uqinisekise ukuba xa sikutsalela umnxeba, uvumela kuphela udliwano-ndlebe ukuba luqale xa ukwindawo yabucala apho uziva ukhululekile ukuthetha ngaphandle kokumanyelwa okanye ukuphazanyiswa. Ukuba ngelixa wenziwa udliwano-ndlebe, uye waphazamiseka, ndicela ucele umphandi ukuba ame ude uzive ukhuselekile ukuqhubeka nokuthetha.\"\n  },\n  {\n    \"old\": \"We would like to know more about your experience with the chatbot. Siyabona ukuba akhange uyigqibe inkqubo kwaye singathanda ukwazi yintoni engesiyenze ngokwahlukileyo ukuphucula amava akho, kwaye siphucule namava abanye abazali abafana nawe kwixesha elizayo. \",\n    \"new\": \"Singathanda ukuva ngakumbi malunga namava wakho ne-chatbot. Siyabona ukuba akhange uyigqibe inkqubo kwaye singathanda ukwazi yintoni engesiyenze ngokwahlukileyo ukuphucula amava akho, kwaye siphucule namava abanye abazali abafana nawe kwixesha elizayo. \"\n  },\n  {\n    \"old\": \" Do I have to agree to be interviewed?\",\n    \"new\": \" Ingaba kufuneka ndivume kudliwano-ndlebe?\"\n  },\n  {\n    \"old\": \"With your permission, we will record the interview to help us remember the discussion and later write down what was said. Siza kucima nayiphi na ingcaciso yobuqu esiyiqokelele kuwe ekupheleni kophononongo kwaye, emva kokubhala udliwano-ndlebe lwakho, sitshintshe nayiphi na idatha enokukhokelela ekukuchazeni kwindawo yokukhuphela. Sinokusebenzisa i-software ye-Artificial Intelligence (AI), iMicrosoft Transcriber, ukukhuphela udliwano-ndlebe ekuqaleni, emva koko siya kujonga/sijongisise oku kukhutshelweyo. Olu lwazi luveliswe yi-AI luya kuqwalaselwa kwaye lugcinwe ngokukhuselekileyo kwiiseva zeDyunivesithi yaseKapa ezikhuselwe ngokuyimfihlo, kwaye ngokungqinelana nePOPIA. Ngamalungu eqela lophando kuphela agunyazisiweyo aya kukwazi ukufikelela kuyo, kwaye le datha iya kuba yeye Global Parenting Initiative kwiDyunivesithi yaseKapa.\",\n    \"new\": \"Ngemvume yakho, sizakurekhoda oludliwano-ndlebe ukusinceda sikhumbule ebesixoxe ngako kwaye kamva sikubhale phantsi obekuthethiwe. Siza kucima nayiphi na ingcaciso yobuqu esiyiqokelele kuwe ekupheleni kophononongo kwaye, emva kokubhala udliwano-ndlebe lwakho, sitshintshe nayiphi na idatha enokukhokelela ekukuchazeni kwindawo yokukhuphela. Sinokusebenzisa i-software ye-Artificial Intelligence (AI), iMicrosoft Transcriber, ukukhuphela udliwano-ndlebe ekuqaleni, emva koko siya kujonga/sijongisise oku kukhutshelweyo. Olu lwazi luveliswe yi-AI luya kuqwalaselwa kwaye lugcinwe ngokukhuselekileyo kwiiseva zeDyunivesithi yaseKapa ezikhuselwe ngokuyimfihlo, kwaye ngokungqinelana nePOPIA. Ngamalungu eqela lophando kuphela agunyazisiweyo aya kukwazi ukufikelela kuyo, kwaye le datha iya kuba yeye Global Parenting Initiative kwiDyunivesithi yaseKapa.\"\n  },\n  {\n    \"old\": \"Siqokelela kuphela oko sikudingayo koluphononongo kwaye sikugcina ngokukhuselekileyo. Your information, like your consent form and interview recording, and any information you provide via email or WhatsApp, will be kept safe on secure servers at the University of Cape Town. \",\n    \"new\": \"Siqokelela kuphela oko sikudingayo koluphononongo kwaye sikugcina ngokukhuselekileyo. Ulwazi lwakho, ukufana nefomu yakho yemvume kunye norekhodingi yodliwano-ndlebe, nayo nayiphi na ingcaciso oyinikeza nge-imeyile okanye nge-WhatsApp, luya kugcinwa lukhuselekile kwiiseva ezikhuselekileyo kwiDyunivesithi yaseKapa. \"\n  },\n  {\n    \"old\": \"Interview recordings will be deleted after we have written our notes. Nayiphi na inkcukacha echaza wena izakugcinwa bucala kwaye ngabasebenzi abagunyazisiweyo kuphela abanokufikelela kuzo. Yonke idatha iya kugcinwa iminyaka emihlanu emva koluphononongo, kodwa inkcukacha zomntu ziya kususwa xa isifundo siphelile. \",\n    \"new\": \"Iirekhodingi zodliwano-ndlebe ziyakucinywa emva kokuba sibhale phantsi amanqaku ethu. Nayiphi na inkcukacha echaza wena izakugcinwa bucala kwaye ngabasebenzi abagunyazisiweyo kuphela abanokufikelela kuzo. Yonke idatha iya kugcinwa iminyaka emihlanu emva koluphononongo, kodwa inkcukacha zomntu ziya kususwa xa isifundo siphelile. \"\n  },\n  {\n    \"old\": \"Ukuthatha kwakho inxaxheba kunye nento osixelela yona izakusinceda siqondisise singazixhasa njani iintsapho ezifana nezakho. We plan to share the results in reports and at conferences so others can learn from this study too.\",\n    \"new\": \"Ukuthatha kwakho inxaxheba kunye nento osixelela yona izakusinceda siqondisise singazixhasa njani iintsapho ezifana nezakho. Sicwangcisa ukwabelana ngeziphumo kwiingxelo nakwii-nkomfa ukuze nabanye bafunde kolu phononongo.\"\n  },\n  {\n    \"old\": \"The principal investigators of this study are Prof Cathy Ward and Cindee Bruyns and the Co-investigator is Carly Katzef all from the University of Cape Town.\",\n    \"new\": \"Abaphononongi abaziintloko kolu phononongo nguNjinga Cathy Ward no Cindee Bruyns ze Co-investigator ngu Carly Katzef bonke basuka kwiDyunivesithi yaseKapa.\"\n  },\n  {\n    \"old\": \"Are there any risks in being interviewed?   \",\n    \"new\": \"Ingaba ikhona imingcipheko ekwenziweni udliwano-ndlebe?   \"\n  },\n  {\n    \"old\": \"Who pays for the study?\",\n    \"new\": \"Ngubani obhatalela oluphononongo?\"\n  },\n  {\n    \"old\": \"University of Cape Town Centre for Social Science Research \",\n    \"new\": \"Idyunivesithi yaseKapa iZiko lezoPhando lweNzululwazi yeZentlalo \"\n  },\n  {\n    \"old\": \"Human Research Ethics Committee\",\n    \"new\": \"Ikomiti yeeNdlela zokuziphatha zoPhando loLuntu\"\n  },\n  {\n    \"old\": \"I am okay with the interview being recorded. I know the recordings will be used for research.\",\n    \"new\": \"Ndilungile ukuba nodliwano-ndlebe olurekhodiweyo. Ndiyazi irekhodingi izakusetyenziswa kuphando.\"\n  }\n];\n\nfor (const { old, new: newText } of replacements) {\n  const results = context.document.body.search(old, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find expected text: \" + old.substring(0, 60));\n  }\n\n  // Replace every match (expected to be exactly one) with the new text.\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM (PowerShell-style) edit script.\n# Applies the same set of exact text replacements as edit.js, using\n# Range.Find.Execute(...) with Replace:=wdReplaceAll (2) against the whole\n# document content. Each `old` string is the full, verbatim text of an\n# existing run/paragraph, so every replacement is unambiguous (exactly one\n# match) and formatting already on the run/paragraph is left untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = 'Amava akho nenkqubo ye-ParentText abalulekile kuphononongo lwethu. Singathanda ukuva malunga namava akho ngayo, kwaye yintoni mhlawumbi engakhange ikusebenzele wena. This interview is part of a study carried out by researchers from the Universities of Cape Town in South Africa and the University of Oxford in the United Kingdom. '; New = 'Amava akho nenkqubo ye-ParentText abalulekile kuphononongo lwethu. Singathanda ukuva malunga namava akho ngayo, kwaye yintoni mhlawumbi engakhange ikusebenzele wena. Oludliwano-ndlebe luyinxalenye yophononongo olwenziwa ngabaphandi abaphuma kwiDyunivesithi yaseKapa eMzantsi Afrika kunye neDyunivesithi yaseOxford eUnited Kingdom. ' }\n    @{ Old = 'Before you decide if you\u2019d like to be interviewed, it\u2019s important for you to know why we\u2019re doing this research and what participating in it would involve. All the information you might need is explained below, but if you have any questions about your participation or our study, please email the study team at swift@globalparenting.org or message us on WhatsApp at +27 XX XXX XXXX. Silapha ukuzokunceda wena!'; New = 'Ngaphambi kokuba ugqibe ekubeni ungathanda na ukuba nodliwano-ndlebe, kubalulekile ukuba wazi kutheni sisenza oluphando nje kwaye ukuthatha inxalenye kungaquka ntoni. Lonke ulwazi ozakuludinga luchaziwe ngezantsi, kodwa ukuba unayo nayiphi na imibuzo malunga nokuthatha inxaxheba okanye ngophando lwethu, ndicela u-imeyilele iqela lophononongo ku swift@globalparenting.org okanye uthumele umyalezo kuthi ku WhatsApp ku +27 XX XXX XXXX. Silapha ukuzokunceda wena!' }\n    @{ Old = 'We would like to have a telephonic conversation with you which will last a maximum of 15 minutes. Omnye wabaphandi bethu uzakutsalela umnxeba athethe nawe ngexesha elikulungeleyo wena. Akukho zimpendulu zilungileyo okanye ezingalunganga, sifuna nje ukuva amava kunye nemibono yakho nge chatbot. Please make sure that when we call, that you only let the interview start when you are in a private space where you feel comfortable to talk without being overheard or interrupted. Ukuba ngelixa wenziwa udliwano-ndlebe, uye waphazamiseka, ndicela ucele umphandi ukuba ame ude uzive ukhuselekile ukuqhubeka nokuthetha.'; New = 'Singathanda ukuba nencoko ngomnxeba nawe ozakuthatha imizuzu eyi-15 ubude. Omnye wabaphandi bethu uzakutsalela umnxeba athethe nawe ngexesha elikulungeleyo wena. Akukho zimpendulu zilungileyo okanye ezingalunganga, sifuna nje ukuva amava kunye nemibono yakho nge chatbot. Nceda uqinisekise ukuba xa sikutsalela umnxeba, uvumela kuphela udliwano-ndlebe ukuba luqale xa ukwindawo yabucala apho uziva ukhululekile ukuthetha ngaphandle kokumanyelwa okanye ukuphazanyiswa. Ukuba ngelixa wenziwa udliwano-ndlebe, uye waphazamiseka, ndicela ucele umphandi ukuba ame ude uzive ukhuselekile ukuqhubeka nokuthetha.' }\n    @{ Old = 'We would like to know more about your experience with the chatbot. Siyabona ukuba akhange uyigqibe inkqubo kwaye singathanda ukwazi yintoni engesiyenze ngokwahlukileyo ukuphucula amava akho, kwaye siphucule namava abanye abazali abafana nawe kwixesha elizayo. '; New = 'Singathanda ukuva ngakumbi malunga namava wakho ne-chatbot. Siyabona ukuba akhange uyigqibe inkqubo kwaye singathanda ukwazi yintoni engesiyenze ngokwahlukileyo ukuphucula amava akho, kwaye siphucule namava abanye abazali abafana nawe kwixesha elizayo. ' }\n    @{ Old = ' Do I have to agree to be interviewed?'; New = ' Ingaba kufuneka ndivume kudliwano-ndlebe?' }\n    @{ Old = 'With your permission, we will record the interview to help us remember the discussion and later write down what was said. Siza kucima nayiphi na ingcaciso yobuqu esiyiqokelele kuwe ekupheleni kophononongo kwaye, emva kokubhala udliwano-ndlebe lwakho, sitshintshe nayiphi na idatha enokukhokelela ekukuchazeni kwindawo yokukhuphela. Sinokusebenzisa i-software ye-Artificial Intelligence (AI), iMicrosoft Transcriber, ukukhuphela udliwano-ndlebe ekuqaleni, emva koko siya kujonga/sijongisise oku kukhutshelweyo. Olu lwazi luveliswe yi-AI luya kuqwalaselwa kwaye lugcinwe ngokukhuselekileyo kwiiseva zeDyunivesithi yaseKapa ezikhuselwe ngokuyimfihlo, kwaye ngokungqinelana nePOPIA. Ngamalungu eqela lophando kuphela agunyazisiweyo aya kukwazi ukufikelela kuyo, kwaye le datha iya kuba yeye Global Parenting Initiative kwiDyunivesithi yaseKapa.'; New = 'Ngemvume yakho, sizakurekhoda oludliwano-ndlebe ukusinceda sikhumbule ebesixoxe ngako kwaye kamva sikubhale phantsi obekuthethiwe. Siza kucima nayiphi na ingcaciso yobuqu esiyiqokelele kuwe ekupheleni kophononongo kwaye, emva kokubhala udliwano-ndlebe lwakho, sitshintshe nayiphi na idatha enokukhokelela ekukuchazeni kwindawo yokukhuphela. Sinokusebenzisa i-software ye-Artificial Intelligence (AI), iMicrosoft Transcriber, ukukhuphela udliwano-ndlebe ekuqaleni, emva koko siya kujonga/sijongisise oku kukhutshelweyo. Olu lwazi luveliswe yi-AI luya kuqwalaselwa kwaye lugcinwe ngokukhuselekileyo kwiiseva zeDyunivesithi yaseKapa ezikhuselwe ngokuyimfihlo, kwaye ngokungqinelana nePOPIA. Ngamalungu eqela lophando kuphela agunyazisiweyo aya kukwazi ukufikelela kuyo, kwaye le datha iya kuba yeye Global Parenting Initiative kwiDyunivesithi yaseKapa.' }\n    @{ Old = 'Siqokelela kuphela oko sikudingayo koluphononongo kwaye sikugcina ngokukhuselekileyo. Your information, like your consent form and interview recording, and any information you provide via email or WhatsApp, will be kept safe on secure servers at the University of Cape Town. '; New = 'Siqokelela kuphela oko sikudingayo koluphononongo kwaye sikugcina ngokukhuselekileyo. Ulwazi lwakho, ukufana nefomu yakho yemvume kunye norekhodingi yodliwano-ndlebe, nayo nayiphi na ingcaciso oyinikeza nge-imeyile okanye nge-WhatsApp, luya kugcinwa lukhuselekile kwiiseva ezikhuselekileyo kwiDyunivesithi yaseKapa. ' }\n    @{ Old = 'Interview recordings will be deleted after we have written our notes. Nayiphi na inkcukacha echaza wena izakugcinwa bucala kwaye ngabasebenzi abagunyazisiweyo kuphela abanokufikelela kuzo. Yonke idatha iya kugcinwa iminyaka emihlanu emva koluphononongo, kodwa inkcukacha zomntu ziya kususwa xa isifundo siphelile. '; New = 'Iirekhodingi zodliwano-ndlebe ziyakucinywa emva kokuba sibhale phantsi amanqaku ethu. Nayiphi na inkcukacha echaza wena izakugcinwa bucala kwaye ngabasebenzi abagunyazisiweyo kuphela abanokufikelela kuzo. Yonke idatha iya kugcinwa iminyaka emihlanu emva koluphononongo, kodwa inkcukacha zomntu ziya kususwa xa isifundo siphelile. ' }\n    @{ Old = 'Ukuthatha kwakho inxaxheba kunye nento osixelela yona izakusinceda siqondisise singazixhasa njani iintsapho ezifana nezakho. We plan to share the results in reports and at conferences so others can learn from this study too.'; New = 'Ukuthatha kwakho inxaxheba kunye nento osixelela yona izakusinceda siqondisise singazixhasa njani iintsapho ezifana nezakho. Sicwangcisa ukwabelana ngeziphumo kwiingxelo nakwii-nkomfa ukuze nabanye bafunde kolu phononongo.' }\n    @{ Old = 'The principal investigators of this study are Prof Cathy Ward and Cindee Bruyns and the Co-investigator is Carly Katzef all from the University of Cape Town.'; New = 'Abaphononongi abaziintloko kolu phononongo nguNjinga Cathy Ward no Cindee Bruyns ze Co-investigator ngu Carly Katzef bonke basuka kwiDyunivesithi yaseKapa.' }\n    @{ Old = 'Are there any risks in being interviewed?   '; New = 'Ingaba ikhona imingcipheko ekwenziweni udliwano-ndlebe?   ' }\n    @{ Old = 'Who pays for the study?'; New = 'Ngubani obhatalela oluphononongo?' }\n    @{ Old = 'University of Cape Town Centre for Social Science Research '; New = 'Idyunivesithi yaseKapa iZiko lezoPhando lweNzululwazi yeZentlalo ' }\n    @{ Old = 'Human Research Ethics Committee'; New = 'Ikomiti yeeNdlela zokuziphatha zoPhando loLuntu' }\n    @{ Old = 'I am okay with the interview being recorded. I know the recordings will be used for research.'; New = 'Ndilungile ukuba nodliwano-ndlebe olurekhodiweyo. Ndiyazi irekhodingi izakusetyenziswa kuphando.' }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $found = $range.Find.Execute(\n        $pair.Old,   # FindText\n        $true,       # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $pair.New,   # ReplaceWith\n        2            # Replace (wdReplaceAll)\n    )\n    if (-not $found) {\n        throw \"Could not find expected text: \" + $pair.Old.Substring(0, [Math]::Min(60, $pair.Old.Length))\n    }\n}\n\n"}
